# The underlying "new TPM" data refresh adds a third sending/target cluster
# ("ECs") to the Efna5 -> Epha3 ligand-receptor pair table, turning the old
# 2x3 sending/target-cluster grid (FAPs, MuSCs) into a full 3x3 grid
# (ECs, FAPs, MuSCs), and refreshes every numeric column for all pairs.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Efna5"
$ws.Cells.Item(2, 3).Value = "Epha3"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.05800433333333333
$ws.Cells.Item(2, 8).Value = 0.174013
$ws.Cells.Item(2, 9).Value = 0.02087975181349295
$ws.Cells.Item(2, 10).Value = 0.02087975181349295
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.003058333333333333
$ws.Cells.Item(2, 14).Value = 0.009175
$ws.Cells.Item(2, 15).Value = 0.0001379486413073712
$ws.Cells.Item(2, 16).Value = 0.0001379486413073712
$ws.Cells.Item(2, 17).Value = 0.0001773965861111111
$ws.Cells.Item(2, 18).Value = 0.001596569275
$ws.Cells.Item(2, 19).Value = 0.000002880333393506473
$ws.Cells.Item(2, 20).Value = 0.000002880333393506473

# Row 3: ECs -> FAPs
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Efna5"
$ws.Cells.Item(3, 3).Value = "Epha3"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.05800433333333333
$ws.Cells.Item(3, 8).Value = 0.174013
$ws.Cells.Item(3, 9).Value = 0.02087975181349295
$ws.Cells.Item(3, 10).Value = 0.02087975181349295
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 22.03620333333333
$ws.Cells.Item(3, 14).Value = 66.10861
$ws.Cells.Item(3, 15).Value = 0.9939610820947024
$ws.Cells.Item(3, 16).Value = 0.9939610820947024
$ws.Cells.Item(3, 17).Value = 1.278195283547778
$ws.Cells.Item(3, 18).Value = 11.50375755193
$ws.Cells.Item(3, 19).Value = 0.02075366070640828
$ws.Cells.Item(3, 20).Value = 0.02075366070640828

# Row 4: ECs -> MuSCs
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Efna5"
$ws.Cells.Item(4, 3).Value = "Epha3"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.05800433333333333
$ws.Cells.Item(4, 8).Value = 0.174013
$ws.Cells.Item(4, 9).Value = 0.02087975181349295
$ws.Cells.Item(4, 10).Value = 0.02087975181349295
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 0.6666666666666666
$ws.Cells.Item(4, 13).Value = 0.130825
$ws.Cells.Item(4, 14).Value = 0.392475
$ws.Cells.Item(4, 15).Value = 0.005900969263990248
$ws.Cells.Item(4, 16).Value = 0.005900969263990248
$ws.Cells.Item(4, 17).Value = 0.007588416908333333
$ws.Cells.Item(4, 18).Value = 0.068295752175
$ws.Cells.Item(4, 19).Value = 0.0001232107736911666
$ws.Cells.Item(4, 20).Value = 0.0001232107736911666

# Row 5: FAPs -> ECs
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Efna5"
$ws.Cells.Item(5, 3).Value = "Epha3"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 1.666083666666667
$ws.Cells.Item(5, 8).Value = 4.998251
$ws.Cells.Item(5, 9).Value = 0.5997381826733804
$ws.Cells.Item(5, 10).Value = 0.5997381826733805
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.003058333333333333
$ws.Cells.Item(5, 14).Value = 0.009175
$ws.Cells.Item(5, 15).Value = 0.0001379486413073712
$ws.Cells.Item(5, 16).Value = 0.0001379486413073712
$ws.Cells.Item(5, 17).Value = 0.005095439213888889
$ws.Cells.Item(5, 18).Value = 0.045858952925
$ws.Cells.Item(5, 19).Value = 0.00008273306743994485
$ws.Cells.Item(5, 20).Value = 0.00008273306743994486

# Row 6: FAPs -> FAPs
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Efna5"
$ws.Cells.Item(6, 3).Value = "Epha3"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 1.666083666666667
$ws.Cells.Item(6, 8).Value = 4.998251
$ws.Cells.Item(6, 9).Value = 0.5997381826733804
$ws.Cells.Item(6, 10).Value = 0.5997381826733805
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 22.03620333333333
$ws.Cells.Item(6, 14).Value = 66.10861
$ws.Cells.Item(6, 15).Value = 0.9939610820947024
$ws.Cells.Item(6, 16).Value = 0.9939610820947024
$ws.Cells.Item(6, 17).Value = 36.71415844901222
$ws.Cells.Item(6, 18).Value = 330.42742604111
$ws.Cells.Item(6, 19).Value = 0.5961164130235435
$ws.Cells.Item(6, 20).Value = 0.5961164130235436

# Row 7: FAPs -> MuSCs
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Efna5"
$ws.Cells.Item(7, 3).Value = "Epha3"
$ws.Cells.Item(7, 4).Value = "MuSCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 1.666083666666667
$ws.Cells.Item(7, 8).Value = 4.998251
$ws.Cells.Item(7, 9).Value = 0.5997381826733804
$ws.Cells.Item(7, 10).Value = 0.5997381826733805
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 0.130825
$ws.Cells.Item(7, 14).Value = 0.392475
$ws.Cells.Item(7, 15).Value = 0.005900969263990248
$ws.Cells.Item(7, 16).Value = 0.005900969263990248
$ws.Cells.Item(7, 17).Value = 0.2179653956916666
$ws.Cells.Item(7, 18).Value = 1.961688561225
$ws.Cells.Item(7, 19).Value = 0.003539036582396986
$ws.Cells.Item(7, 20).Value = 0.003539036582396987

# Row 8: MuSCs -> ECs
$ws.Cells.Item(8, 1).Value = "MuSCs"
$ws.Cells.Item(8, 2).Value = "Efna5"
$ws.Cells.Item(8, 3).Value = "Epha3"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 1.053930333333333
$ws.Cells.Item(8, 8).Value = 3.161791
$ws.Cells.Item(8, 9).Value = 0.3793820655131266
$ws.Cells.Item(8, 10).Value = 0.3793820655131266
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.003058333333333333
$ws.Cells.Item(8, 14).Value = 0.009175
$ws.Cells.Item(8, 15).Value = 0.0001379486413073712
$ws.Cells.Item(8, 16).Value = 0.0001379486413073712
$ws.Cells.Item(8, 17).Value = 0.003223270269444445
$ws.Cells.Item(8, 18).Value = 0.029009432425
$ws.Cells.Item(8, 19).Value = 0.00005233524047391992
$ws.Cells.Item(8, 20).Value = 0.00005233524047391992

# Row 9: MuSCs -> FAPs
$ws.Cells.Item(9, 1).Value = "MuSCs"
$ws.Cells.Item(9, 2).Value = "Efna5"
$ws.Cells.Item(9, 3).Value = "Epha3"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 1.053930333333333
$ws.Cells.Item(9, 8).Value = 3.161791
$ws.Cells.Item(9, 9).Value = 0.3793820655131266
$ws.Cells.Item(9, 10).Value = 0.3793820655131266
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 22.03620333333333
$ws.Cells.Item(9, 14).Value = 66.10861
$ws.Cells.Item(9, 15).Value = 0.9939610820947024
$ws.Cells.Item(9, 16).Value = 0.9939610820947024
$ws.Cells.Item(9, 17).Value = 23.22462312450111
$ws.Cells.Item(9, 18).Value = 209.02160812051
$ws.Cells.Item(9, 19).Value = 0.3770910083647506
$ws.Cells.Item(9, 20).Value = 0.3770910083647506

# Row 10: MuSCs -> MuSCs
$ws.Cells.Item(10, 1).Value = "MuSCs"
$ws.Cells.Item(10, 2).Value = "Efna5"
$ws.Cells.Item(10, 3).Value = "Epha3"
$ws.Cells.Item(10, 4).Value = "MuSCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 1.053930333333333
$ws.Cells.Item(10, 8).Value = 3.161791
$ws.Cells.Item(10, 9).Value = 0.3793820655131266
$ws.Cells.Item(10, 10).Value = 0.3793820655131266
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 0.6666666666666666
$ws.Cells.Item(10, 13).Value = 0.130825
$ws.Cells.Item(10, 14).Value = 0.392475
$ws.Cells.Item(10, 15).Value = 0.005900969263990248
$ws.Cells.Item(10, 16).Value = 0.005900969263990248
$ws.Cells.Item(10, 17).Value = 0.1378804358583333
$ws.Cells.Item(10, 18).Value = 1.240923922725
$ws.Cells.Item(10, 19).Value = 0.002238721907902095
$ws.Cells.Item(10, 20).Value = 0.002238721907902095
